{"js": "// Locate the \"RNF02\" text in the requirement list and split it after\n// \"RN\", placing the insertion point (tracked by Word's \"_GoBack\"\n// bookmark) between \"RN\" and \"F02\" - mirroring what Word does when the\n// last edit in the document happened inside that run.\nconst body = context.document.body;\nconst matches = body.search(\"RNF02\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length > 0) {\n  const target = matches.items[0];\n\n  const prefixMatches = target.search(\"RN\", { matchCase: true });\n  prefixMatches.load(\"items\");\n  await context.sync();\n\n  const prefix = prefixMatches.items[0];\n  const splitPoint = prefix.getRange(\"After\");\n\n  // Word keeps a single \"_GoBack\" bookmark marking the most recent edit\n  // location; move it here (removing it from wherever it was before).\n  context.document.deleteBookmark(\"_GoBack\");\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"RNF02\" text in the requirement list and split it after \"RN\",\n# placing the insertion point (tracked by Word's \"_GoBack\" bookmark)\n# between \"RN\" and \"F02\" - mirroring what Word does when the last edit in\n# the document happened inside that run.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"RNF02\")\n\nif ($rng.Find.Found) {\n    $splitPos = $rng.Start + 2\n    $insertRange = $d.Range($splitPos, $splitPos)\n\n    # Word keeps a single \"_GoBack\" bookmark marking the most recent edit\n    # location; move it here (removing it from wherever it was before).\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks(\"_GoBack\").Delete()\n    }\n    $d.Bookmarks.Add(\"_GoBack\", $insertRange)\n}\n"}
